$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) contain text values that can look like
# plain numbers (e.g. "324.66") or percentages. Excel would normally try to
# be "smart" and coerce such literals into floating point numbers, which
# loses formatting (trailing zeros, thousand-dot groups, etc.). Temporarily
# force the range to Text format so every assignment below is stored as a
# literal string, matching the original inline-string cell type.
$numRange = $ws.Range("D2:E51")
$numRange.NumberFormat = "@"

$ws.Range('D2').Value = '27.730.77'
$ws.Range('E2').Value = '  +3.06%  '
$ws.Range('D3').Value = '1.866.70'
$ws.Range('E3').Value = '  +3.13%  '
$ws.Range('E4').Value = '  +2.83%  '
$ws.Range('D5').Value = '324.66'
$ws.Range('E5').Value = '  +4.39%  '
$ws.Range('E6').Value = '  +2.64%  '
$ws.Range('D7').Value = '0.4415'
$ws.Range('E7').Value = '  +2.93%  '
$ws.Range('D8').Value = '0.3804'
$ws.Range('E8').Value = '  +3.21%  '
$ws.Range('D9').Value = '0.07459'
$ws.Range('E9').Value = '  +3.08%  '
$ws.Range('D10').Value = '0.8847'
$ws.Range('E10').Value = '  +2.45%  '
$ws.Range('D11').Value = '21.83'
$ws.Range('D12').Value = '1.883.60'
$ws.Range('E12').Value = '  -8.37%  '
$ws.Range('D13').Value = '5.559'
$ws.Range('E13').Value = '  +3.06%  '
$ws.Range('D14').Value = '6.756'
$ws.Range('E14').Value = '  +2.07%  '
$ws.Range('D15').Value = '0.07207'
$ws.Range('E15').Value = '  +3.47%  '
$ws.Range('D16').Value = '83.94'
$ws.Range('E16').Value = '  +3.85%  '
$ws.Range('E17').Value = '  +2.55%  '
$ws.Range('D18').Value = '0.000009099'
$ws.Range('E18').Value = '  +2.06%  '
$ws.Range('D19').Value = '1.034'
$ws.Range('E19').Value = '  +2.71%  '
$ws.Range('D20').Value = '15.52'
$ws.Range('E20').Value = '  +2.34%  '
$ws.Range('D21').Value = '27.765.37'
$ws.Range('E21').Value = '  +3.02%  '
$ws.Range('D22').Value = '5.323'
$ws.Range('E22').Value = '  +2.57%  '
$ws.Range('D23').Value = '11.45'
$ws.Range('E23').Value = '  +4.83%  '
$ws.Range('D24').Value = '158.09'
$ws.Range('E24').Value = '  +2.59%  '
$ws.Range('D25').Value = '1.947'
$ws.Range('E25').Value = '  +3.19%  '
$ws.Range('D26').Value = '18.85'
$ws.Range('E26').Value = '  +2.64%  '
$ws.Range('E27').Value = '  +3.97%  '
$ws.Range('D28').Value = '5.329'
$ws.Range('E28').Value = '  +1.88%  '
$ws.Range('D29').Value = '117.69'
$ws.Range('E29').Value = '  +2.61%  '
$ws.Range('D30').Value = '0.09102'
$ws.Range('E30').Value = '  +1.69%  '
$ws.Range('D31').Value = '1.216'
$ws.Range('E31').Value = '  +5.03%  '
$ws.Range('D32').Value = '0.7702'
$ws.Range('E32').Value = '  +3.81%  '
$ws.Range('D33').Value = '3.008'
$ws.Range('E33').Value = '  +7.07%  '
$ws.Range('D34').Value = '4.582'
$ws.Range('E34').Value = '  +3.49%  '
$ws.Range('D35').Value = '1.035'
$ws.Range('E35').Value = '  +2.71%  '
$ws.Range('D36').Value = '1.164'
$ws.Range('E36').Value = '  +3.72%  '
$ws.Range('D37').Value = '0.01993'
$ws.Range('E37').Value = '  +3.78%  '
$ws.Range('D38').Value = '0.05357'
$ws.Range('E38').Value = '  +2.45%  '
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').Value = '0.5200'
$ws.Range('E39').Value = '  +2.30%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '2.841'
$ws.Range('E40').Value = '  +3.38%  '
$ws.Range('D41').Value = '0.1695'
$ws.Range('E41').Value = '  +2.61%  '
$ws.Range('D42').Value = '6.851'
$ws.Range('E42').Value = '  +5.86%  '
$ws.Range('D43').Value = '8.706'
$ws.Range('E43').Value = '  +4.85%  '
$ws.Range('D44').Value = '109.89'
$ws.Range('E44').Value = '  +2.30%  '
$ws.Range('D45').Value = '10.66'
$ws.Range('E45').Value = '  +2.64%  '
$ws.Range('D46').Value = '1.730'
$ws.Range('E46').Value = '  +5.21%  '
$ws.Range('D47').Value = '0.4701'
$ws.Range('E47').Value = '  +2.84%  '
$ws.Range('D48').Value = '0.06430'
$ws.Range('E48').Value = '  +2.47%  '
$ws.Range('D49').Value = '1.881'
$ws.Range('E49').Value = '  +3.69%  '
$ws.Range('D50').Value = '39.78'
$ws.Range('E50').Value = '  +4.60%  '
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').Value = '0.9365'
$ws.Range('E51').Value = '  +2.68%  '

# Put the cell style back to the workbook default so the cells are not left
# with an explicit "Text" style that was not present in the original file.
$numRange.Style = "Normal"

